# New crime data collected - weekly CompStat update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich text cells: replace only the changed characters) ---
# A8: "Volume 30   Number  15" -> "...16"
$ws.Range("A8").Characters(21, 2).Text = "16"

# C9: "Report Covering the Week  4/10/2023  Through  4/16/2023"
#  -> "Report Covering the Week  4/17/2023  Through  4/23/2023"
$ws.Range("C9").Characters(27, 9).Text = "4/17/2023"
$ws.Range("C9").Characters(47, 9).Text = "4/23/2023"

# --- Data table updates (rows 14-29) ---

# Row 14 (Murder)
# C14 switches from a numeric count to the literal text "0" (style matches
# the existing text-"0" cells elsewhere in the table, e.g. D22).
$ws.Range("D22").Copy($ws.Range("C14"))
$ws.Range("F14").Value = 2
$ws.Range("N14").Value = -16.666666666666

# Row 15 (Rape)
$ws.Range("F15").Value = 2
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -30

# Row 16 (Robbery)
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 500
$ws.Range("F16").Value = 15
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = 87.5
$ws.Range("I16").Value = 65
$ws.Range("J16").Value = 29
$ws.Range("K16").Value = 124.137931034483
$ws.Range("L16").Value = 195.454545454545
$ws.Range("M16").Value = -18.75
$ws.Range("N16").Value = -79.938271604938

# Row 17 (Fel. Assault)
$ws.Range("C17").Value = 12
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 300
$ws.Range("F17").Value = 33
$ws.Range("G17").Value = 18
$ws.Range("H17").Value = 83.333333333333
$ws.Range("I17").Value = 118
$ws.Range("J17").Value = 99
$ws.Range("K17").Value = 19.191919191919
$ws.Range("L17").Value = 61.643835616438
$ws.Range("M17").Value = 20.408163265306
$ws.Range("N17").Value = -43.809523809523

# Row 18 (Burglary)
$ws.Range("D22").Copy($ws.Range("C18"))
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 14
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = 55.555555555555
$ws.Range("I18").Value = 52
$ws.Range("J18").Value = 27
$ws.Range("K18").Value = 92.592592592592
$ws.Range("L18").Value = 160
$ws.Range("M18").Value = -27.777777777777
$ws.Range("N18").Value = -90.076335877862

# Row 19 (Gr. Larceny)
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 14.285714285714
$ws.Range("F19").Value = 29
$ws.Range("G19").Value = 19
$ws.Range("H19").Value = 52.631578947368
$ws.Range("I19").Value = 127
$ws.Range("J19").Value = 110
$ws.Range("K19").Value = 15.454545454545
$ws.Range("L19").Value = 115.254237288136
$ws.Range("M19").Value = 15.454545454545
$ws.Range("N19").Value = -19.620253164557

# Row 20 (G.L.A.)
$ws.Range("C20").Value = 3
$ws.Range("E20").Value = 50
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = 12.5
$ws.Range("I20").Value = 36
$ws.Range("J20").Value = 24
$ws.Range("L20").Value = 63.636363636363
$ws.Range("M20").Value = -30.769230769230
$ws.Range("N20").Value = -91.705069124424

# Row 21 (TOTAL)
$ws.Range("C21").Value = 29
$ws.Range("D21").Value = 16
$ws.Range("E21").Value = 81.25
$ws.Range("F21").Value = 104
$ws.Range("G21").Value = 62
$ws.Range("H21").Value = 67.741935483871
$ws.Range("I21").Value = 410
$ws.Range("J21").Value = 292
$ws.Range("K21").Value = 40.410958904109
$ws.Range("L21").Value = 99.029126213592
$ws.Range("M21").Value = -3.301886792452
$ws.Range("N21").Value = -75.566150178784

# Row 23 (Transit)
$ws.Range("C23").Value = 2
# D23/E23 switch from literal text ("0" / "***.*") to real numbers; borrow
# the number style from existing numeric cells in the same columns first.
$ws.Range("G16").Copy($ws.Range("D23"))
$ws.Range("D23").Value = 4
$ws.Range("K14").Copy($ws.Range("E23"))
$ws.Range("E23").Value = -50
$ws.Range("F23").Value = 5
$ws.Range("G23").Value = 5
$ws.Range("I23").Value = 29
$ws.Range("J23").Value = 19
$ws.Range("K23").Value = 52.631578947368
$ws.Range("L23").Value = 70.588235294117
$ws.Range("M23").Value = 141.666666666667

# Row 24 (Housing)
$ws.Range("C24").Value = 28
$ws.Range("D24").Value = 22
$ws.Range("E24").Value = 27.272727272727
$ws.Range("F24").Value = 90
$ws.Range("G24").Value = 86
$ws.Range("H24").Value = 4.651162790697
$ws.Range("I24").Value = 369
$ws.Range("J24").Value = 312
$ws.Range("K24").Value = 18.269230769230
$ws.Range("L24").Value = 63.274336283185
$ws.Range("M24").Value = -6.106870229007

# Row 25 (Petit Larceny)
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = 25
$ws.Range("F25").Value = 62
$ws.Range("G25").Value = 57
$ws.Range("H25").Value = 8.771929824561
$ws.Range("I25").Value = 216
$ws.Range("J25").Value = 196
$ws.Range("K25").Value = 10.204081632653
$ws.Range("L25").Value = 68.75
$ws.Range("M25").Value = -40.166204986149

# Row 26 (Misd. Assault)
$ws.Range("F26").Value = 2
$ws.Range("H26").Value = 100
$ws.Range("L26").Value = -10

# Row 27 (UCR Rape*)
$ws.Range("C27").Value = 2
$ws.Range("E27").Value = 100
$ws.Range("I27").Value = 21
$ws.Range("J27").Value = 18
$ws.Range("K27").Value = 16.666666666666
$ws.Range("L27").Value = 23.529411764705

# Row 28 (Other Sex Crimes)
$ws.Range("D22").Copy($ws.Range("C28"))
$ws.Range("F28").Value = 2
$ws.Range("L28").Value = -60
$ws.Range("N28").Value = -78.947368421052

# Row 29 (Shooting Vic.)
$ws.Range("D22").Copy($ws.Range("C29"))
$ws.Range("F29").Value = 2
$ws.Range("L29").Value = -60
$ws.Range("N29").Value = -76.470588235294
